$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 200
$ws.Range("I33").Value = 200
$ws.Range("K33").Value = 200
$ws.Range("M33").Value = 29
$ws.Range("H76").Value = 3455.625
$ws.Range("I76").Value = 3416.25
$ws.Range("J76").Value = 3495
$ws.Range("K76").Value = 3416.25
$ws.Range("L76").Value = 3495
$ws.Range("M76").Value = -3101.25
$ws.Range("N76").Value = -4125
$ws.Range("H79").Value = 3455.625
$ws.Range("I79").Value = 3416.25
$ws.Range("J79").Value = 3495
$ws.Range("K79").Value = 3416.25
$ws.Range("L79").Value = 3495
$ws.Range("M79").Value = -2324.25
$ws.Range("N79").Value = -5679
$ws.Range("H98").Value = 709.1875
$ws.Range("I98").Value = 742.0769
$ws.Range("K98").Value = 742.0769
$ws.Range("M98").Value = 755.9231
$ws.Range("H122").Value = 709.1875
$ws.Range("I122").Value = 742.0769
$ws.Range("K122").Value = 2226.2307
$ws.Range("M122").Value = 223.7692999999999
$ws.Range("H129").Value = 1037.1045
$ws.Range("I129").Value = 440.44446
$ws.Range("J129").Value = 1129.6897
$ws.Range("K129").Value = 1321.33338
$ws.Range("L129").Value = 3389.0691
$ws.Range("M129").Value = 3678.66662
$ws.Range("N129").Value = -13389.0691
$ws.Range("H137").Value = 1545.0358
$ws.Range("I137").Value = 1365.0416
$ws.Range("J137").Value = 2625
$ws.Range("K137").Value = 4095.1248
$ws.Range("L137").Value = 7875
$ws.Range("M137").Value = -1545.1248
$ws.Range("N137").Value = -12975
$ws.Range("H138").Value = 2392.8281
$ws.Range("I138").Value = 2292.6428
$ws.Range("J138").Value = 2420.88
$ws.Range("K138").Value = 6877.928400000001
$ws.Range("L138").Value = 7262.64
$ws.Range("M138").Value = -1737.928400000001
$ws.Range("N138").Value = -17542.64

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1503.4615
$ws.Range("I2").Value = 1388.8788
$ws.Range("K2").Value = 1388.8788
$ws.Range("M2").Value = -1275.8788
$ws.Range("H32").Value = 6798.711
$ws.Range("I32").Value = 4917.256
$ws.Range("K32").Value = 4917.256
$ws.Range("M32").Value = -4630.256
$ws.Range("H45").Value = 3988.6
$ws.Range("I45").Value = 5252.1665
$ws.Range("J45").Value = 3146.2222
$ws.Range("K45").Value = 5252.1665
$ws.Range("L45").Value = 3146.2222
$ws.Range("M45").Value = -4875.1665
$ws.Range("N45").Value = -3900.2222
$ws.Range("H74").Value = 47621824
$ws.Range("I74").Value = 76926744
$ws.Range("K74").Value = 76926744
$ws.Range("M74").Value = -76925870
$ws.Range("H77").Value = 47621824
$ws.Range("I77").Value = 76926744
$ws.Range("K77").Value = 384633720
$ws.Range("M77").Value = -384629352
$ws.Range("H102").Value = 800
$ws.Range("I102").Value = 800
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 800
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 822
$ws.Range("N102").ClearContents()
$ws.Range("H105").Value = 38031.285
$ws.Range("J105").Value = 38031.285
$ws.Range("L105").Value = 38031.285
$ws.Range("N105").Value = -45019.285
$ws.Range("H116").Value = 1503.4615
$ws.Range("I116").Value = 1388.8788
$ws.Range("K116").Value = 1388.8788
$ws.Range("M116").Value = 905.1212
$ws.Range("H132").Value = 16645.824
$ws.Range("I132").Value = 1833.1034
$ws.Range("J132").Value = 102559.6
$ws.Range("K132").Value = 5499.3102
$ws.Range("L132").Value = 307678.8
$ws.Range("M132").Value = -2969.3102
$ws.Range("N132").Value = -312738.8

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1503.4615
$ws.Range("I3").Value = 1388.8788
$ws.Range("K3").Value = 1388.8788
$ws.Range("M3").Value = -1274.8788
$ws.Range("H20").Value = 3780.182
$ws.Range("I20").Value = 4410.222
$ws.Range("J20").Value = 945
$ws.Range("K20").Value = 4410.222
$ws.Range("L20").Value = 945
$ws.Range("M20").Value = -4163.222
$ws.Range("N20").Value = -1439
$ws.Range("H86").Value = 1722.6666
$ws.Range("I86").Value = 1454.7727
$ws.Range("K86").Value = 1454.7727
$ws.Range("M86").Value = -331.7727
$ws.Range("H89").Value = 1722.6666
$ws.Range("I89").Value = 1454.7727
$ws.Range("K89").Value = 7273.863499999999
$ws.Range("M89").Value = -1657.863499999999

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1333.6364
$ws.Range("I16").Value = 1340
$ws.Range("J16").Value = 1322.5
$ws.Range("K16").Value = 1340
$ws.Range("L16").Value = 1322.5
$ws.Range("M16").Value = -1053
$ws.Range("N16").Value = -1896.5
$ws.Range("H31").Value = 23057.611
$ws.Range("I31").Value = 108600.336
$ws.Range("K31").Value = 108600.336
$ws.Range("M31").Value = -108305.336
$ws.Range("H34").Value = 23057.611
$ws.Range("I34").Value = 108600.336
$ws.Range("K34").Value = 108600.336
$ws.Range("M34").Value = -108398.336
$ws.Range("H62").Value = 76926424
$ws.Range("I62").Value = 83336210
$ws.Range("J62").Value = 9006
$ws.Range("K62").Value = 83336210
$ws.Range("L62").Value = 9006
$ws.Range("M62").Value = -83335586
$ws.Range("N62").Value = -10254
$ws.Range("H65").Value = 76926424
$ws.Range("I65").Value = 83336210
$ws.Range("J65").Value = 9006
$ws.Range("K65").Value = 416681050
$ws.Range("L65").Value = 45030
$ws.Range("M65").Value = -416677930
$ws.Range("N65").Value = -51270
$ws.Range("H99").Value = 25005050
$ws.Range("I99").Value = 4076.8462
$ws.Range("J99").Value = 71435430
$ws.Range("K99").Value = 4076.8462
$ws.Range("L99").Value = 71435430
$ws.Range("M99").Value = -2578.8462
$ws.Range("N99").Value = -71438426
$ws.Range("H113").Value = 1333.6364
$ws.Range("I113").Value = 1340
$ws.Range("J113").Value = 1322.5
$ws.Range("K113").Value = 1340
$ws.Range("L113").Value = 1322.5
$ws.Range("M113").Value = 830
$ws.Range("N113").Value = -5662.5
$ws.Range("H122").Value = 1337.2084
$ws.Range("I122").Value = 1423.625
$ws.Range("K122").Value = 4270.875
$ws.Range("M122").Value = -1820.875
$ws.Range("H126").Value = 25005050
$ws.Range("I126").Value = 4076.8462
$ws.Range("J126").Value = 71435430
$ws.Range("K126").Value = 12230.5386
$ws.Range("L126").Value = 214306290
$ws.Range("M126").Value = -9760.5386
$ws.Range("N126").Value = -214311230
$ws.Range("H132").Value = 14775.22
$ws.Range("I132").Value = 19837.178
$ws.Range("K132").Value = 59511.534
$ws.Range("M132").Value = -56981.534

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7500045
$ws.Range("I4").Value = 90
$ws.Range("K4").Value = 270
$ws.Range("M4").Value = -158
$ws.Range("H122").Value = 452.875
$ws.Range("J122").Value = 969.4286
$ws.Range("L122").Value = 8724.857399999999
$ws.Range("N122").Value = -13624.8574
$ws.Range("H131").Value = 780.86
$ws.Range("I131").Value = 300
$ws.Range("J131").Value = 785.71716
$ws.Range("K131").Value = 900
$ws.Range("L131").Value = 2357.15148
$ws.Range("M131").Value = 4140
$ws.Range("N131").Value = -12437.15148

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1719.9487
$ws.Range("J113").Value = 2021.4286
$ws.Range("L113").Value = 2021.4286
$ws.Range("N113").Value = -6361.4286
$ws.Range("H122").Value = 133336070
$ws.Range("I122").Value = 111112080
$ws.Range("J122").Value = 142860640
$ws.Range("K122").Value = 333336240
$ws.Range("L122").Value = 428581920
$ws.Range("M122").Value = -333333790
$ws.Range("N122").Value = -428586820
$ws.Range("H126").Value = 4435.6763
$ws.Range("I126").Value = 3390.8635
$ws.Range("J126").Value = 6351.1665
$ws.Range("K126").Value = 10172.5905
$ws.Range("L126").Value = 19053.4995
$ws.Range("M126").Value = -7702.5905
$ws.Range("N126").Value = -23993.4995
$ws.Range("H132").Value = 26670.305
$ws.Range("I132").Value = 5356.6113
$ws.Range("J132").Value = 103399.6
$ws.Range("K132").Value = 16069.8339
$ws.Range("L132").Value = 310198.8
$ws.Range("M132").Value = -13539.8339
$ws.Range("N132").Value = -315258.8

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4754.115
$ws.Range("I7").Value = 3379.0715
$ws.Range("J7").Value = 6358.3335
$ws.Range("K7").Value = 3379.0715
$ws.Range("L7").Value = 6358.3335
$ws.Range("M7").Value = -3267.0715
$ws.Range("N7").Value = -6582.3335
$ws.Range("H82").Value = 3767
$ws.Range("J82").Value = 4000.75
$ws.Range("L82").Value = 4000.75
$ws.Range("N82").Value = -4722.75
$ws.Range("H85").Value = 3767
$ws.Range("J85").Value = 4000.75
$ws.Range("L85").Value = 4000.75
$ws.Range("N85").Value = -6496.75
$ws.Range("H122").Value = 936439.8
$ws.Range("I122").Value = 1785057.4
$ws.Range("J122").Value = 2960.5
$ws.Range("K122").Value = 5355172.199999999
$ws.Range("L122").Value = 8881.5
$ws.Range("M122").Value = -5352722.199999999
$ws.Range("N122").Value = -13781.5
$ws.Range("H126").Value = 4754.115
$ws.Range("I126").Value = 3379.0715
$ws.Range("J126").Value = 6358.3335
$ws.Range("K126").Value = 10137.2145
$ws.Range("L126").Value = 19075.0005
$ws.Range("M126").Value = -7667.2145
$ws.Range("N126").Value = -24015.0005

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5051360
$ws.Range("I107").Value = 1374.75
$ws.Range("J107").Value = 9091349
$ws.Range("K107").Value = 4124.25
$ws.Range("L107").Value = 27274047
$ws.Range("M107").Value = -2204.25
$ws.Range("N107").Value = -27277887
$ws.Range("H126").Value = 989.7778
$ws.Range("I126").Value = 985
$ws.Range("J126").Value = 999.3333
$ws.Range("K126").Value = 2955
$ws.Range("L126").Value = 2997.9999
$ws.Range("M126").Value = -485
$ws.Range("N126").Value = -7937.9999
